$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Cells.Item(2, 2).Value = 1.273737299345328
$ws.Cells.Item(2, 3).Value = 0.3063437075418847
$ws.Cells.Item(2, 5).Value = 0.1007883983159239
$ws.Cells.Item(2, 6).Value = 0.4443680307746121
$ws.Cells.Item(2, 7).Value = 0.2485704132827209
$ws.Cells.Item(2, 8).Value = 0.4385839356070704
$ws.Cells.Item(2, 9).Value = 0.3972548585594522
$ws.Cells.Item(2, 12).Value = 0.2154034315090172
$ws.Cells.Item(2, 15).Value = 1.29235046425886

# Row 3
$ws.Cells.Item(3, 2).Value = 1.126921057944401
$ws.Cells.Item(3, 3).Value = 0.2913928432706996
$ws.Cells.Item(3, 5).Value = 0.1022443921992622
$ws.Cells.Item(3, 6).Value = 0.3878228170618172
$ws.Cells.Item(3, 7).Value = 0.2537946314256629
$ws.Cells.Item(3, 8).Value = 0.4452073174105493
$ws.Cells.Item(3, 9).Value = 0.4070851709228798
$ws.Cells.Item(3, 12).Value = 0.2045791096859233
$ws.Cells.Item(3, 15).Value = 1.317043984179094

# Row 4
$ws.Cells.Item(4, 2).Value = 1.036464842328257
$ws.Cells.Item(4, 3).Value = 0.2822237813277582
$ws.Cells.Item(4, 5).Value = 0.1032103372163595
$ws.Cells.Item(4, 6).Value = 0.3531389305169483
$ws.Cells.Item(4, 7).Value = 0.2573350179696341
$ws.Cells.Item(4, 8).Value = 0.4495651765221176
$ws.Cells.Item(4, 9).Value = 0.4134934390738567
$ws.Cells.Item(4, 12).Value = 0.1980181924961073
$ws.Cells.Item(4, 15).Value = 1.333512258986381

# Row 5
$ws.Cells.Item(5, 2).Value = 0.9995279378966302
$ws.Cells.Item(5, 3).Value = 0.2784904288306507
$ws.Cells.Item(5, 5).Value = 0.1036220571482325
$ws.Cells.Item(5, 6).Value = 0.3390132514313251
$ws.Cells.Item(5, 7).Value = 0.2588611184792349
$ws.Cells.Item(5, 8).Value = 0.4514141767436755
$ws.Cells.Item(5, 9).Value = 0.4161983324653917
$ws.Cells.Item(5, 12).Value = 0.1953661542435015
$ws.Cells.Item(5, 15).Value = 1.340550960283302

# Row 6
$ws.Cells.Item(6, 2).Value = 0.9933901303133439
$ws.Cells.Item(6, 3).Value = 0.2778707081949676
$ws.Cells.Item(6, 5).Value = 0.103691515336843
$ws.Cells.Item(6, 6).Value = 0.336668177824194
$ws.Cells.Item(6, 7).Value = 0.2591195526327219
$ws.Cells.Item(6, 8).Value = 0.4517256174804025
$ws.Cells.Item(6, 9).Value = 0.4166531182943691
$ws.Cells.Item(6, 12).Value = 0.1949270936350018
$ws.Cells.Item(6, 15).Value = 1.341739504309707

# Row 7
$ws.Cells.Item(7, 2).Value = 1.035966999321147
$ws.Cells.Item(7, 3).Value = 0.2821734188991059
$ws.Cells.Item(7, 5).Value = 0.1032158165784267
$ws.Cells.Item(7, 6).Value = 0.3529483938344953
$ws.Cells.Item(7, 7).Value = 0.2573552623476125
$ws.Cells.Item(7, 8).Value = 0.4495898167301746
$ws.Cells.Item(7, 9).Value = 0.4135295399795034
$ws.Cells.Item(7, 12).Value = 0.197982338570057
$ws.Cells.Item(7, 15).Value = 1.333605859337887

# Row 8
$ws.Cells.Item(8, 2).Value = 1.223181123815095
$ws.Cells.Item(8, 3).Value = 0.3011866310298785
$ws.Cells.Item(8, 5).Value = 0.1012754869596808
$ws.Cells.Item(8, 6).Value = 0.4248636149813478
$ws.Cells.Item(8, 7).Value = 0.2503024963225329
$ws.Cells.Item(8, 8).Value = 0.4408072317359171
$ws.Cells.Item(8, 9).Value = 0.4005669513034755
$ws.Cells.Item(8, 12).Value = 0.2116535953146439
$ws.Cells.Item(8, 15).Value = 1.300593209673607

# Row 9
$ws.Cells.Item(9, 2).Value = 1.587739031613694
$ws.Cells.Item(9, 3).Value = 0.3385421786347251
$ws.Cells.Item(9, 5).Value = 0.09804168386543921
$ws.Cells.Item(9, 6).Value = 0.5661985755041457
$ws.Cells.Item(9, 7).Value = 0.2391249625452474
$ws.Cells.Item(9, 8).Value = 0.4258965288243743
$ws.Cells.Item(9, 9).Value = 0.3781097380250849
$ws.Cells.Item(9, 12).Value = 0.2391343705540976
$ws.Cells.Item(9, 15).Value = 1.246253858284007

# Row 10
$ws.Cells.Item(10, 2).Value = 1.85389911457969
$ws.Cells.Item(10, 3).Value = 0.3660128330928387
$ws.Cells.Item(10, 5).Value = 0.09601428965141778
$ws.Cells.Item(10, 6).Value = 0.6702781546542269
$ws.Cells.Item(10, 7).Value = 0.2325480498444321
$ws.Cells.Item(10, 8).Value = 0.4163540710123499
$ws.Cells.Item(10, 9).Value = 0.3634259661681591
$ws.Cells.Item(10, 12).Value = 0.2597293837862509
$ws.Cells.Item(10, 15).Value = 1.212714544750114

# Row 11
$ws.Cells.Item(11, 2).Value = 1.974594882679469
$ws.Cells.Item(11, 3).Value = 0.3785119954132767
$ws.Cells.Item(11, 5).Value = 0.09516770265140018
$ws.Cells.Item(11, 6).Value = 0.7176906081379002
$ws.Cells.Item(11, 7).Value = 0.2299149288767453
$ws.Cells.Item(11, 8).Value = 0.4123202976567768
$ws.Cells.Item(11, 9).Value = 0.3571421905571972
$ws.Cells.Item(11, 12).Value = 0.2691857057761524
$ws.Cells.Item(11, 15).Value = 1.198852325934581

# Row 12
$ws.Cells.Item(12, 2).Value = 2.020241906434649
$ws.Cells.Item(12, 3).Value = 0.3832451328644026
$ws.Cells.Item(12, 5).Value = 0.09485801151140016
$ws.Cells.Item(12, 6).Value = 0.7356546913071611
$ws.Cells.Item(12, 7).Value = 0.2289697339516863
$ws.Cells.Item(12, 8).Value = 0.4108370400106836
$ws.Cells.Item(12, 9).Value = 0.3548198080863605
$ws.Cells.Item(12, 12).Value = 0.2727790405239148
$ws.Cells.Item(12, 15).Value = 1.193804455654657

# Row 13
$ws.Cells.Item(13, 2).Value = 2.010413622381463
$ws.Cells.Item(13, 3).Value = 0.3822257744467095
$ws.Cells.Item(13, 5).Value = 0.09492422442019155
$ws.Cells.Item(13, 6).Value = 0.7317853510981394
$ws.Cells.Item(13, 7).Value = 0.229170984720092
$ws.Cells.Item(13, 8).Value = 0.4111545170372324
$ws.Cells.Item(13, 9).Value = 0.3553174297327937
$ws.Cells.Item(13, 12).Value = 0.2720046012454986
$ws.Cells.Item(13, 15).Value = 1.194882631698633

# Row 14
$ws.Cells.Item(14, 2).Value = 1.978351466852985
$ws.Cells.Item(14, 3).Value = 0.3789013962475565
$ws.Cells.Item(14, 5).Value = 0.09514200580895249
$ws.Cells.Item(14, 6).Value = 0.7191683204515869
$ws.Cells.Item(14, 7).Value = 0.2298361249075214
$ws.Cells.Item(14, 8).Value = 0.4121973819770091
$ws.Cells.Item(14, 9).Value = 0.3569499802398104
$ws.Cells.Item(14, 12).Value = 0.2694810833615264
$ws.Cells.Item(14, 15).Value = 1.19843299235373

# Row 15
$ws.Cells.Item(15, 2).Value = 1.95870485085328
$ws.Cells.Item(15, 3).Value = 0.3768651036336337
$ws.Cells.Item(15, 5).Value = 0.09527682206778643
$ws.Cells.Item(15, 6).Value = 0.7114413442032514
$ws.Cells.Item(15, 7).Value = 0.230250312354741
$ws.Cells.Item(15, 8).Value = 0.4128419311678755
$ws.Cells.Item(15, 9).Value = 0.3579574130230316
$ws.Cells.Item(15, 12).Value = 0.2679369699619656
$ws.Cells.Item(15, 15).Value = 1.200633950916384

# Row 16
$ws.Cells.Item(16, 2).Value = 1.846003427547942
$ws.Cells.Item(16, 3).Value = 0.3651960018387967
$ws.Cells.Item(16, 5).Value = 0.09607113989934923
$ws.Cells.Item(16, 6).Value = 0.6671810134426437
$ws.Cells.Item(16, 7).Value = 0.2327273771179534
$ws.Cells.Item(16, 8).Value = 0.4166238747308952
$ws.Cells.Item(16, 9).Value = 0.3638446111286093
$ws.Cells.Item(16, 12).Value = 0.2591131390341275
$ws.Cells.Item(16, 15).Value = 1.213648617124377

# Row 17
$ws.Cells.Item(17, 2).Value = 1.776764867442068
$ws.Cells.Item(17, 3).Value = 0.3580377710554501
$ws.Cells.Item(17, 5).Value = 0.09657781974153679
$ws.Cells.Item(17, 6).Value = 0.6400460337125793
$ws.Cells.Item(17, 7).Value = 0.2343390910482341
$ws.Cells.Item(17, 8).Value = 0.4190226962337533
$ws.Cells.Item(17, 9).Value = 0.3675577776407764
$ws.Cells.Item(17, 12).Value = 0.2537223147122631
$ws.Cells.Item(17, 15).Value = 1.221990591457967

# Row 18
$ws.Cells.Item(18, 2).Value = 1.736904917587708
$ws.Cells.Item(18, 3).Value = 0.3539208167417485
$ws.Cells.Item(18, 5).Value = 0.09687637183464837
$ws.Cells.Item(18, 6).Value = 0.6244449056556647
$ws.Cells.Item(18, 7).Value = 0.2352998543032641
$ws.Cells.Item(18, 8).Value = 0.42043133957403
$ws.Cells.Item(18, 9).Value = 0.3697307485637795
$ws.Cells.Item(18, 12).Value = 0.250629904941178
$ws.Cells.Item(18, 15).Value = 1.22691992614017

# Row 19
$ws.Cells.Item(19, 2).Value = 1.723402971497023
$ws.Cells.Item(19, 3).Value = 0.3525269455399211
$ws.Cells.Item(19, 5).Value = 0.09697867957287087
$ws.Cells.Item(19, 6).Value = 0.619163680173358
$ws.Cells.Item(19, 7).Value = 0.2356309386116067
$ws.Cells.Item(19, 8).Value = 0.4209132440069965
$ws.Cells.Item(19, 9).Value = 0.3704728714719767
$ws.Cells.Item(19, 12).Value = 0.2495842898681104
$ws.Cells.Item(19, 15).Value = 1.228611430654936

# Row 20
$ws.Cells.Item(20, 2).Value = 1.784139152930152
$ws.Cells.Item(20, 3).Value = 0.3587997518426675
$ws.Cells.Item(20, 5).Value = 0.09652314550775998
$ws.Cells.Item(20, 6).Value = 0.642933953830422
$ws.Cells.Item(20, 7).Value = 0.2341640262072602
$ws.Cells.Item(20, 8).Value = 0.4187643454945444
$ws.Cells.Item(20, 9).Value = 0.3671586473691821
$ws.Cells.Item(20, 12).Value = 0.2542953244903003
$ws.Cells.Item(20, 15).Value = 1.22108898378724

# Row 21
$ws.Cells.Item(21, 2).Value = 1.987770492341326
$ws.Cells.Item(21, 3).Value = 0.3798778505409643
$ws.Cells.Item(21, 5).Value = 0.09507774242786127
$ws.Cells.Item(21, 6).Value = 0.7228739723491628
$ws.Cells.Item(21, 7).Value = 0.2296393457017558
$ws.Cells.Item(21, 8).Value = 0.4118898658495382
$ws.Cells.Item(21, 9).Value = 0.3564689079555778
$ws.Cells.Item(21, 12).Value = 0.2702219656256659
$ws.Cells.Item(21, 15).Value = 1.197384691259742

# Row 22
$ws.Cells.Item(22, 2).Value = 2.120516835203205
$ws.Cells.Item(22, 3).Value = 0.3936532828370503
$ws.Cells.Item(22, 5).Value = 0.09419658493401073
$ws.Cells.Item(22, 6).Value = 0.7751780083420101
$ws.Cells.Item(22, 7).Value = 0.2269849299630877
$ws.Cells.Item(22, 8).Value = 0.4076549246474173
$ws.Cells.Item(22, 9).Value = 0.3498157388402028
$ws.Cells.Item(22, 12).Value = 0.280703329521657
$ws.Cells.Item(22, 15).Value = 1.183067169839134

# Row 23
$ws.Cells.Item(23, 2).Value = 2.049699553585981
$ws.Cells.Item(23, 3).Value = 0.3863012365823977
$ws.Cells.Item(23, 5).Value = 0.09466106216769887
$ws.Cells.Item(23, 6).Value = 0.7472568307830727
$ws.Cells.Item(23, 7).Value = 0.2283738360806922
$ws.Cells.Item(23, 8).Value = 0.4098915661028215
$ws.Cells.Item(23, 9).Value = 0.3533361038034597
$ws.Cells.Item(23, 12).Value = 0.2751026574996729
$ws.Cells.Item(23, 15).Value = 1.190600943254211

# Row 24
$ws.Cells.Item(24, 2).Value = 1.780805405754222
$ws.Cells.Item(24, 3).Value = 0.3584552652581863
$ws.Cells.Item(24, 5).Value = 0.09654784113080339
$ws.Cells.Item(24, 6).Value = 0.6416283278902171
$ws.Cells.Item(24, 7).Value = 0.2342430666177862
$ws.Cells.Item(24, 8).Value = 0.4188810538890309
$ws.Cells.Item(24, 9).Value = 0.3673389751220384
$ws.Cells.Item(24, 12).Value = 0.2540362454316067
$ws.Cells.Item(24, 15).Value = 1.22149618507882

# Row 25
$ws.Cells.Item(25, 2).Value = 1.489403392767713
$ws.Cells.Item(25, 3).Value = 0.3284307758003706
$ws.Cells.Item(25, 5).Value = 0.09885533569605798
$ws.Cells.Item(25, 6).Value = 0.5279251897347166
$ws.Cells.Item(25, 7).Value = 0.2418629077081889
$ws.Cells.Item(25, 8).Value = 0.4296824480512527
$ws.Cells.Item(25, 9).Value = 0.383866872308654
$ws.Cells.Item(25, 12).Value = 0.2316286694593117
$ws.Cells.Item(25, 15).Value = 1.259836239717316
